$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number-format on the Price/Volume columns so that
# numeric-looking strings (e.g. "212.98") are stored as text, matching
# the original inline-string cell contents instead of being coerced to numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.297.37'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.65%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.608.98'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.35%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.98'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.15%  '
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.05%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.25%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.78%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.49'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +2.62%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0814'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.26%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.832.52'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.31%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.627.19'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.40%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.03'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.53%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.517'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.19%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '26.265.84'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.51%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.28'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +3.10%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0₃0727'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.82%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.00%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '201.83'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.16%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.27'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +1.10%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.34'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.61%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.64%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.44%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '143.26'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.21%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.03%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.122'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.10%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.25'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.61%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +2.36%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0497'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +5.45%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.01%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.20'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +2.93%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.85%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.40%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +1.52%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.162.42'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +3.28%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +2.21%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.01%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.93%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.71%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.91%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.38'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +4.17%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.13%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.743.82'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.22%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '92.17'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.51%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +14.30%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +1.30%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '54.05'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.13%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.48%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.25%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.14%  '
